$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# Update existing row 2 values (future production -> uncon_planned_qty / con_planned_qty / produced_qty)
$ws.Range("G2").Value = 710
$ws.Range("H2").Value = 710
$ws.Range("J2").Value = 675

# Add new row 3 for MAT_B / LINE_B
$ws.Range("A3").Value = "MAT_B"
$ws.Range("B3").Value = "PLANT_001"
$ws.Range("C3").Value = "LINE_B"

$ws.Range("D3").Value = 45295
$ws.Range("E3").Value = 45296
$ws.Range("F3").Value = 45297

$ws.Range("D3:F3").NumberFormat = $ws.Range("D2:F2").NumberFormat

$ws.Range("G3").Value = 104
$ws.Range("H3").Value = 104
$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = 92
